$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Variante 1 / Variante 2" example block to the right of the
# existing "Service Bytes" table (columns M:O, rows 16-19).
# Values are entered in the same order the original author typed them so
# that the generated shared-string table indices line up with the target.
$ws.Range("M17").Value = " "
$ws.Range("N16").Value = "Beispiel"
$ws.Range("O17").Value = "0x8912345678"
$ws.Range("O18").Value = "0x89 0x12345678"
$ws.Range("N18").Value = "Variante 2:"
$ws.Range("N17").Value = "Variante 1:"
$ws.Range("N19").Value = "Var1 oder Var2?"

# Widen the two new columns to fit the new content.
$ws.Columns.Item(14).ColumnWidth = 13.6
$ws.Columns.Item(15).ColumnWidth = 14.6

# Restore the view state (scrolled down/right a bit, new active selection).
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("M23").Select()
